$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Drop the oldest quarter column (D: "6 months ended 1399/06") and shift
#        everything one column to the left (E..M -> D..L). ---
$ws.Columns("D").Delete()

# --- 2. Re-create column M (now blank) with the same per-row formatting as
#        column L, so the newly appended quarter matches the existing style. ---
$ws.Range("L1:L28").Copy()
$ws.Range("M1:M28").PasteSpecial(-4122)

# --- 3. Header row: add the label for the new right-most quarter. ---
$ws.Range("M8").Value = "12 ماهه منتهی به 1401/12"

# --- 4. Publish-date row: update the restated values and add the new one. ---
$ws.Range("K9").Value = "1402-02-20 (4)"
$ws.Range("L9").Value = "1402-02-29 (9)"
$ws.Range("M9").Value = "1402-02-29"

# --- 5. Data rows: column K carries a restated figure (recalculated under the
#        new read_price algorithm) and column M carries the brand-new quarter. ---
$ws.Range("K11").Value = 2417138
$ws.Range("M11").Value = 5033548

$ws.Range("K12").Value = -2303126
$ws.Range("M12").Value = -4073861

$ws.Range("K13").Value = 114012
$ws.Range("M13").Value = 959687

$ws.Range("K14").Value = -138345
$ws.Range("M14").Value = -342452

$ws.Range("K15").Value = 0
$ws.Range("M15").Value = 0

$ws.Range("K16").Value = 77248
$ws.Range("M16").Value = -379117

$ws.Range("K17").Value = 52915
$ws.Range("M17").Value = 238118

$ws.Range("K18").Value = -64647
$ws.Range("M18").Value = -330348

$ws.Range("K19").Value = 12244
$ws.Range("M19").Value = 113977

$ws.Range("K20").Value = 512
$ws.Range("M20").Value = 21747

$ws.Range("K21").Value = 0
$ws.Range("M21").Value = -3914

$ws.Range("K22").Value = 512
$ws.Range("M22").Value = 17833

$ws.Range("K23").Value = 0
$ws.Range("M23").Value = 0

$ws.Range("K24").Value = 512
$ws.Range("M24").Value = 17833

# Row 25 (EPS): G25/K25 become literal "-" placeholders, M25 is the new figure.
$ws.Range("G25").Value = "-"
$ws.Range("K25").Value = "-"
$ws.Range("M25").Value = 14

$ws.Range("K26").Value = 0
$ws.Range("M26").Value = 1273000

$ws.Range("K27").Value = 0
$ws.Range("M27").Value = 14
